$wb = $excel.ActiveWorkbook

# --- Sheet 1: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("Q15").Value = 431.64
$ws1.Range("Q22").Value = "1 de 20"

# --- Sheet 2: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F15").Value = 431.64
$ws2.Range("F22").Value = 3840.97

# --- Sheet 3: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D14").Value = 431.64
$ws3.Range("E14").Value = 51.36000000000001
$ws3.Range("F14").Value = 0.8936645962732919
$ws3.Range("D19").Value = 3840.970000000001
$ws3.Range("E19").Value = 46546.22762291769
$ws3.Range("F19").Value = 0.07622908558528381
